$wb = $excel.ActiveWorkbook

# --- Sheet 1: addVisitor -------------------------------------------------
# Replace the old "record/col count" columns (B/C) with two rows of
# recently-added visitor e-mail addresses in column A.
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").Value = "brandNewVisitorEmailId"
$ws1.Range("A2").Value = "SeleniumTest+v20191017153943@gmail.com"
$ws1.Range("A3").Value = "SeleniumTest+v20191017154305@gmail.com"

# Clear out the cells that used to hold sqlRecordCount/sqlColCount data.
$ws1.Range("B1:C2").Clear()

# --- Sheet 2: sqlCount (new) ---------------------------------------------
# Add the new sheet right after "addVisitor" and populate it with the
# sqlRecordCount / sqlColCount values that used to live on sheet 1.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "sqlCount"

$newSheet.Range("A1").Value = "sqlRecordCount"
$newSheet.Range("B1").Value = "sqlColCount"
# Leading apostrophe forces these numeric-looking values to be stored as
# text (shared strings), matching the source data (which came from a
# text-typed SQL count, not a numeric cell).
$newSheet.Range("A2").Value = "'252"
$newSheet.Range("B2").Value = "'5"

$ws1.Activate()
